# Scheduled-runner update: refresh computed profit/cost figures on a
# handful of Leve rows across several job sheets (ALC, ARM, BSM, CRP,
# GSM, LTW, WVR) to reflect current market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2153
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10228

$ws.Range("H20").Value = 4666.3335
$ws.Range("I20").Value = 4666.3335
$ws.Range("K20").Value = 4666.3335
$ws.Range("M20").Value = -4436.3335

$ws.Range("H35").Value = 4666.3335
$ws.Range("I35").Value = 4666.3335
$ws.Range("K35").Value = 4666.3335
$ws.Range("M35").Value = -4287.3335

$ws.Range("H40").Value = 4153.4546
$ws.Range("J40").Value = 6549.75
$ws.Range("L40").Value = 6549.75
$ws.Range("N40").Value = -6899.75

$ws.Range("H98").Value = 696.725
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H107").Value = 1002.2727
$ws.Range("J107").Value = 222
$ws.Range("L107").Value = 222
$ws.Range("N107").Value = -4062

$ws.Range("H115").Value = 1040
$ws.Range("I115").Value = 380
$ws.Range("K115").Value = 1140
$ws.Range("M115").Value = 427

$ws.Range("H122").Value = 696.725
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H127").Value = 1193.8
$ws.Range("I127").Value = 1193.8
$ws.Range("K127").Value = 3581.4
$ws.Range("M127").Value = 1378.6

$ws.Range("H129").Value = 1275.5
$ws.Range("I129").Value = 1275.5
$ws.Range("K129").Value = 3826.5
$ws.Range("M129").Value = 1173.5

$ws.Range("H132").Value = 6881.8
$ws.Range("I132").Value = 1712
$ws.Range("K132").Value = 5136
$ws.Range("M132").Value = -2606

$ws.Range("H137").Value = 1638.65
$ws.Range("I137").Value = 1332.6875
$ws.Range("K137").Value = 3998.0625
$ws.Range("M137").Value = -1448.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 60000
$ws.Range("J24").Value = 60000
$ws.Range("L24").Value = 60000
$ws.Range("N24").Value = -60748

$ws.Range("H46").Value = 6717
$ws.Range("I46").Value = 3693.6
$ws.Range("J46").Value = 9740.4
$ws.Range("K46").Value = 3693.6
$ws.Range("L46").Value = 9740.4
$ws.Range("M46").Value = -3374.6
$ws.Range("N46").Value = -10378.4

$ws.Range("H50").Value = 527.8
$ws.Range("J50").Value = 463
$ws.Range("L50").Value = 463
$ws.Range("N50").Value = -1891

$ws.Range("H74").Value = 4135.3335
$ws.Range("I74").Value = 2856.1724
$ws.Range("K74").Value = 2856.1724
$ws.Range("M74").Value = -1982.1724

$ws.Range("H77").Value = 4135.3335
$ws.Range("I77").Value = 2856.1724
$ws.Range("K77").Value = 14280.862
$ws.Range("M77").Value = -9912.861999999999

$ws.Range("H97").Value = 1147.1
$ws.Range("I97").Value = 731.44446
$ws.Range("K97").Value = 731.44446
$ws.Range("M97").Value = -235.44446

$ws.Range("H100").Value = 60000
$ws.Range("J100").Value = 60000
$ws.Range("L100").Value = 60000
$ws.Range("N100").Value = -62164

$ws.Range("H132").Value = 3560.88
$ws.Range("I132").Value = 2602.4783
$ws.Range("K132").Value = 7807.4349
$ws.Range("M132").Value = -5277.4349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 716.3158
$ws.Range("I22").Value = 729.2143
$ws.Range("K22").Value = 729.2143
$ws.Range("M22").Value = -556.2143

$ws.Range("H86").Value = 3034.2666
$ws.Range("I86").Value = 2429
$ws.Range("J86").Value = 6968.5
$ws.Range("K86").Value = 2429
$ws.Range("L86").Value = 6968.5
$ws.Range("M86").Value = -1306
$ws.Range("N86").Value = -9214.5

$ws.Range("H89").Value = 3034.2666
$ws.Range("I89").Value = 2429
$ws.Range("J89").Value = 6968.5
$ws.Range("K89").Value = 12145
$ws.Range("L89").Value = 34842.5
$ws.Range("M89").Value = -6529
$ws.Range("N89").Value = -46074.5

$ws.Range("H107").Value = 6833.3335
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 325
$ws.Range("I19").Value = 325
$ws.Range("K19").Value = 325
$ws.Range("M19").Value = -155

$ws.Range("H24").Value = 325
$ws.Range("I24").Value = 325
$ws.Range("K24").Value = 325
$ws.Range("M24").Value = -155

$ws.Range("H105").Value = 2983
$ws.Range("I105").Value = 2599.75
$ws.Range("J105").Value = 3749.5
$ws.Range("K105").Value = 2599.75
$ws.Range("L105").Value = 3749.5
$ws.Range("M105").Value = -852.75
$ws.Range("N105").Value = -7243.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 145987860
$ws.Range("J11").Value = 353750
$ws.Range("L11").Value = 353750
$ws.Range("N11").Value = -354028

$ws.Range("H80").Value = 3708.65
$ws.Range("I80").Value = 3343.111
$ws.Range("K80").Value = 3343.111
$ws.Range("M80").Value = -2345.111

$ws.Range("H83").Value = 3708.65
$ws.Range("I83").Value = 3343.111
$ws.Range("K83").Value = 16715.555
$ws.Range("M83").Value = -11723.555

$ws.Range("H113").Value = 801921
$ws.Range("I113").Value = 801921
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 801921
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -799751
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1937.4445
$ws.Range("I132").Value = 1316.1177
$ws.Range("K132").Value = 3948.3531
$ws.Range("M132").Value = -1418.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3433
$ws.Range("I7").Value = 3433
$ws.Range("K7").Value = 3433
$ws.Range("M7").Value = -3321

$ws.Range("H40").Value = 3180.9092
$ws.Range("I40").Value = 2979.1
$ws.Range("K40").Value = 2979.1
$ws.Range("M40").Value = -2843.1

$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50386

$ws.Range("H46").Value = 13074.185
$ws.Range("I46").Value = 8260.333000000001
$ws.Range("J46").Value = 17406.65
$ws.Range("K46").Value = 8260.333000000001
$ws.Range("L46").Value = 17406.65
$ws.Range("M46").Value = -8072.333000000001
$ws.Range("N46").Value = -17782.65

$ws.Range("H68").Value = 2200
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2200
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H100").Value = 7583.3335
$ws.Range("J100").Value = 7500
$ws.Range("L100").Value = 7500
$ws.Range("N100").Value = -8582

$ws.Range("H126").Value = 3433
$ws.Range("I126").Value = 3433
$ws.Range("K126").Value = 10299
$ws.Range("M126").Value = -7829

$ws.Range("H132").Value = 13665
$ws.Range("I132").Value = 14422.682
$ws.Range("K132").Value = 43268.046
$ws.Range("M132").Value = -40738.046

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 31597.4
$ws.Range("I45").Value = 7999
$ws.Range("J45").Value = 37497
$ws.Range("K45").Value = 7999
$ws.Range("L45").Value = 37497
$ws.Range("M45").Value = -7508
$ws.Range("N45").Value = -38479

$ws.Range("H107").Value = 1127.8948
$ws.Range("I107").Value = 619
$ws.Range("J107").Value = 1498
$ws.Range("K107").Value = 1857
$ws.Range("L107").Value = 4494
$ws.Range("M107").Value = 63
$ws.Range("N107").Value = -8334

$ws.Range("H122").Value = 5037.1816
$ws.Range("I122").Value = 5094.5
$ws.Range("K122").Value = 15283.5
$ws.Range("M122").Value = -12833.5
